$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.224535333333333
$ws.Range("H2").Value = 12.673606
$ws.Range("I2").Value = 0.7043225486309714
$ws.Range("J2").Value = 0.7043225486309715
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 33.380049
$ws.Range("N2").Value = 100.140147
$ws.Range("O2").Value = 0.3891462059670435
$ws.Range("P2").Value = 0.3891462059670435
$ws.Range("Q2").Value = 141.015196428898
$ws.Range("R2").Value = 1269.136767860082
$ws.Range("S2").Value = 0.274084447576781
$ws.Range("T2").Value = 0.2740844475767811
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.224535333333333
$ws.Range("H3").Value = 12.673606
$ws.Range("I3").Value = 0.7043225486309714
$ws.Range("J3").Value = 0.7043225486309715
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.85786133333334
$ws.Range("N3").Value = 101.573584
$ws.Range("O3").Value = 0.3947165649764305
$ws.Range("P3").Value = 0.3947165649764305
$ws.Range("Q3").Value = 143.0337315137671
$ws.Range("R3").Value = 1287.303583623904
$ws.Range("S3").Value = 0.278007777031062
$ws.Range("T3").Value = 0.278007777031062
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.224535333333333
$ws.Range("H4").Value = 12.673606
$ws.Range("I4").Value = 0.7043225486309714
$ws.Range("J4").Value = 0.7043225486309715
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.53974466666667
$ws.Range("N4").Value = 55.61923400000001
$ws.Range("O4").Value = 0.2161372290565261
$ws.Range("P4").Value = 0.2161372290565261
$ws.Range("Q4").Value = 78.32180641531156
$ws.Range("R4").Value = 704.896257737804
$ws.Range("S4").Value = 0.1522303240231285
$ws.Range("T4").Value = 0.1522303240231285
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.773477
$ws.Range("H5").Value = 5.320431
$ws.Range("I5").Value = 0.2956774513690286
$ws.Range("J5").Value = 0.2956774513690286
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.380049
$ws.Range("N5").Value = 100.140147
$ws.Range("O5").Value = 0.3891462059670435
$ws.Range("P5").Value = 0.3891462059670435
$ws.Range("Q5").Value = 59.198749160373
$ws.Range("R5").Value = 532.788742443357
$ws.Range("S5").Value = 0.1150617583902625
$ws.Range("T5").Value = 0.1150617583902625
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.773477
$ws.Range("H6").Value = 5.320431
$ws.Range("I6").Value = 0.2956774513690286
$ws.Range("J6").Value = 0.2956774513690286
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 33.85786133333334
$ws.Range("N6").Value = 101.573584
$ws.Range("O6").Value = 0.3947165649764305
$ws.Range("P6").Value = 0.3947165649764305
$ws.Range("Q6").Value = 60.04613834385601
$ws.Range("R6").Value = 540.4152450947041
$ws.Range("S6").Value = 0.1167087879453685
$ws.Range("T6").Value = 0.1167087879453685
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.773477
$ws.Range("H7").Value = 5.320431
$ws.Range("I7").Value = 0.2956774513690286
$ws.Range("J7").Value = 0.2956774513690286
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.53974466666667
$ws.Range("N7").Value = 55.61923400000001
$ws.Range("O7").Value = 0.2161372290565261
$ws.Range("P7").Value = 0.2161372290565261
$ws.Range("Q7").Value = 32.879810752206
$ws.Range("R7").Value = 295.9182967698541
$ws.Range("S7").Value = 0.06390690503339758
$ws.Range("T7").Value = 0.06390690503339758
